$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the initial estimate value for "That User Story" (row 3) - now blank
$ws.Range("C3").Value = $null

# Add estimates for time on profile creation (rows 19-23)
$ws.Range("C19").Value = 6
$ws.Range("C20").Value = 6
$ws.Range("C21").Value = 6
$ws.Range("C22").Value = 4
$ws.Range("C23").Value = 4

# Update selection to reflect new active cell
$ws.Range("C22").Select()
